$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 168, shifting existing rows 168-247 down to 169-248
$ws.Rows.Item(168).EntireRow.Insert()

# Populate the newly inserted row 168 with the new data point
$ws.Cells.Item(168, 1).Value = 4
$ws.Cells.Item(168, 2).Value = "Feria Lagunitas de Puerto Montt"
$ws.Cells.Item(168, 3).Value = "Los Lagos"
$ws.Cells.Item(168, 4).Value = 44489
$ws.Cells.Item(168, 5).Value = 10
$ws.Cells.Item(168, 6).Value = 100114001
$ws.Cells.Item(168, 7).Value = "Papa"
$ws.Cells.Item(168, 8).Value = "Asterix"
$ws.Cells.Item(168, 9).Value = "1a (guarda)"
$ws.Cells.Item(168, 10).Value = 100
$ws.Cells.Item(168, 11).Value = 9000
$ws.Cells.Item(168, 12).Value = 9000
$ws.Cells.Item(168, 13).Value = 9000
$ws.Cells.Item(168, 14).Value = "$/saco 25 kilos"
$ws.Cells.Item(168, 15).Value = "Provincia de Llanquihue"
$ws.Cells.Item(168, 16).Value = 360
$ws.Cells.Item(168, 17).Value = 25
$ws.Cells.Item(168, 18).Value = "Hortaliza"
